# Merge Josh's heatmap + Tim's documentation changes into the Scheduled
# Tasks workbook.
#
# 1. Rows 20-24 ("Completed" column E) get a "Yes" flag, copying the
#    existing style used on the "Yes" cells above (E19 etc.) so the
#    shared-style index matches (s="1").
# 2. Five brand-new tasks (rows 25-29) are appended with the same
#    per-column formatting already used by the rows above them (plain
#    for Task ID / Task Name / Assigned Team Member, date style copied
#    from the Deadline column above for column D).
# 3. The worksheet Table (Table_1) is resized from A1:E24 to A1:E29 so
#    the new rows are included.
# 4. The view is left selecting C27, matching where the edit left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Flag existing rows 20-24 as Completed = Yes -------------------
$completedRows = 20, 21, 22, 23, 24
foreach ($r in $completedRows) {
    $cell = $ws.Cells.Item($r, 5)           # column E
    $cell.Value = "Yes"
}
# Copy the style from the "Yes" cell directly above the block (E19) so
# the new cells pick up the same cellXf (s="1") as the rest of the
# column instead of the default style.
$ws.Range("E19").Copy() | Out-Null
$ws.Range("E20:E24").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Append the five new tasks (rows 25-29) ------------------------
$newTasks = @(
    @{ Row = 25; Id = 24; Name = "Work on Documenation, Complete BSL induction"; Member = "Timothy Finn" },
    @{ Row = 26; Id = 25; Name = "Work on Java implementation, Complete BSL induction"; Member = "Yeqing Liu" },
    @{ Row = 27; Id = 26; Name = "Work on LED implementation, Complete BSL induction"; Member = "Josh Francis" },
    @{ Row = 28; Id = 27; Name = "Work on Raspberry Pi interface, Complete BSL induction"; Member = "Michael Douglas" },
    @{ Row = 29; Id = 28; Name = "Work on brain diagram, Complete BSL induction"; Member = "George Proios" }
)

foreach ($task in $newTasks) {
    $r = $task.Row
    $ws.Cells.Item($r, 1).Value = $task.Id
    $ws.Cells.Item($r, 2).Value = $task.Name
    $ws.Cells.Item($r, 3).Value = $task.Member
    $ws.Cells.Item($r, 4).Value = 43349
}

# Copy the date format from the Deadline cell directly above (D24) onto
# the newly written Deadline cells so they carry the same s="7" style.
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D25:D29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Grow the table to cover the new rows --------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E29"))

# --- 4. Leave the selection where the author left it ------------------
$ws.Range("C27").Select()
